$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artículos")

# Row 4 is currently a partially-filled record for "Leche ultrapasteurizada entera Manfrey".
# Fill in the remaining attributes for this article.
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "lt."
$ws.Range("H4").Value = "Sachet"
$ws.Range("I4").Value = "Leches"
$ws.Range("J4").Value = "Argentina"
$ws.Range("K4").Value = 12
$ws.Range("M4").Value = $true
$ws.Range("N4").Value = "C:\VentaSoft\Imágenes de artículos\7791058000595.png"
$ws.Range("O4").Value = $true
